$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column A with the row index values (1-5) for rows 2-6
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Move the active selection from B6 to A7, as in the final saved state
$ws.Range("A7").Select() | Out-Null
